$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new student record on row 2
$ws.Range("A2").Value = 96
$ws.Range("B2").Value = "Багин Максим Евгеньевич"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = "8А"

# Row 2 grows taller to fit the wrapped full name
$ws.Rows.Item(2).RowHeight = 30.6

# Move the active selection to E10
$ws.Range("E10").Select()
